$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H43").Value = 11111711
$ws.Range("I43").Value = 749.75
$ws.Range("J43").Value = 55555556
$ws.Range("K43").Value = 749.75
$ws.Range("L43").Value = 55555556
$ws.Range("M43").Value = -680.75
$ws.Range("N43").Value = -55555694

$ws.Range("H111").Value = 2819.7727
$ws.Range("I111").Value = 1276
$ws.Range("K111").Value = 3828
$ws.Range("M111").Value = -761

$ws.Range("H132").Value = 5488.1924
$ws.Range("I132").Value = 2699.6667
$ws.Range("J132").Value = 17200
$ws.Range("K132").Value = 8099.000100000001
$ws.Range("L132").Value = 51600
$ws.Range("M132").Value = -5569.000100000001
$ws.Range("N132").Value = -56660

$ws.Range("H138").Value = 2103.3333
$ws.Range("J138").Value = 2174.3647
$ws.Range("L138").Value = 6523.0941
$ws.Range("N138").Value = -16803.0941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4253.25
$ws.Range("I32").Value = 4253.25
$ws.Range("K32").Value = 4253.25
$ws.Range("M32").Value = -3966.25

$ws.Range("H45").Value = 1853.875
$ws.Range("I45").Value = 1761.5714
$ws.Range("K45").Value = 1761.5714
$ws.Range("M45").Value = -1384.5714

$ws.Range("H61").Value = 1049.1666
$ws.Range("I61").Value = 828.0714
$ws.Range("K61").Value = 828.0714
$ws.Range("M61").Value = -616.0714

$ws.Range("H102").Value = 27795430
$ws.Range("I102").Value = 33354196
$ws.Range("K102").Value = 33354196
$ws.Range("M102").Value = -33352574

$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H136").Value = 1049.1666
$ws.Range("I136").Value = 828.0714
$ws.Range("K136").Value = 2484.2142
$ws.Range("M136").Value = 65.78579999999965

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 144273120
$ws.Range("I105").Value = 144273120
$ws.Range("K105").Value = 144273120
$ws.Range("M105").Value = -144271373

$ws.Range("H134").Value = 5053.815
$ws.Range("I134").Value = 1182.2858
$ws.Range("J134").Value = 18604.166
$ws.Range("K134").Value = 3546.8574
$ws.Range("L134").Value = 55812.49800000001
$ws.Range("M134").Value = -1011.8574
$ws.Range("N134").Value = -60882.49800000001

$ws.Range("H138").Value = 85390
$ws.Range("J138").Value = 85390
$ws.Range("L138").Value = 85390
$ws.Range("N138").Value = -95670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1319.9459
$ws.Range("I31").Value = 1035.3334
$ws.Range("J31").Value = 1514
$ws.Range("K31").Value = 1035.3334
$ws.Range("L31").Value = 1514
$ws.Range("M31").Value = -740.3334
$ws.Range("N31").Value = -2104

$ws.Range("H34").Value = 1319.9459
$ws.Range("I34").Value = 1035.3334
$ws.Range("J34").Value = 1514
$ws.Range("K34").Value = 1035.3334
$ws.Range("L34").Value = 1514
$ws.Range("M34").Value = -833.3334
$ws.Range("N34").Value = -1918

$ws.Range("H58").Value = 996.55554
$ws.Range("I58").Value = 1101
$ws.Range("J58").Value = 819
$ws.Range("K58").Value = 1101
$ws.Range("L58").Value = 819
$ws.Range("M58").Value = -898
$ws.Range("N58").Value = -1225

$ws.Range("H86").Value = 2787910.5
$ws.Range("I86").Value = 4168433.8
$ws.Range("J86").Value = 26864.375
$ws.Range("K86").Value = 4168433.8
$ws.Range("L86").Value = 26864.375
$ws.Range("M86").Value = -4167310.8
$ws.Range("N86").Value = -29110.375

$ws.Range("H89").Value = 2787910.5
$ws.Range("I89").Value = 4168433.8
$ws.Range("J89").Value = 26864.375
$ws.Range("K89").Value = 20842169
$ws.Range("L89").Value = 134321.875
$ws.Range("M89").Value = -20836553
$ws.Range("N89").Value = -145553.875

$ws.Range("H122").Value = 999
$ws.Range("I122").Value = 998
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2994
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -544
$ws.Range("N122").Value = -7900

$ws.Range("H136").Value = 996.55554
$ws.Range("I136").Value = 1101
$ws.Range("J136").Value = 819
$ws.Range("K136").Value = 3303
$ws.Range("L136").Value = 2457
$ws.Range("M136").Value = -753
$ws.Range("N136").Value = -7557

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 550.16
$ws.Range("J113").Value = 596.44446
$ws.Range("L113").Value = 1789.33338
$ws.Range("N113").Value = -6129.33338

$ws.Range("H132").Value = 1400
$ws.Range("J132").Value = 1400
$ws.Range("L132").Value = 12600
$ws.Range("N132").Value = -17660

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3831644.2
$ws.Range("I7").Value = 5928571.5
$ws.Range("J7").Value = 569757.25
$ws.Range("K7").Value = 5928571.5
$ws.Range("L7").Value = 569757.25
$ws.Range("M7").Value = -5928459.5
$ws.Range("N7").Value = -569981.25

$ws.Range("H8").Value = 3831644.2
$ws.Range("I8").Value = 5928571.5
$ws.Range("J8").Value = 569757.25
$ws.Range("K8").Value = 5928571.5
$ws.Range("L8").Value = 569757.25
$ws.Range("M8").Value = -5928432.5
$ws.Range("N8").Value = -570035.25

$ws.Range("H102").Value = 1758.95
$ws.Range("I102").Value = 1805
$ws.Range("K102").Value = 1805
$ws.Range("M102").Value = -183

$ws.Range("H122").Value = 3106.0527
$ws.Range("I122").Value = 2001.875
$ws.Range("K122").Value = 6005.625
$ws.Range("M122").Value = -3555.625

$ws.Range("H126").Value = 1860.3889
$ws.Range("I126").Value = 1619.9
$ws.Range("J126").Value = 2161
$ws.Range("K126").Value = 4859.700000000001
$ws.Range("L126").Value = 6483
$ws.Range("M126").Value = -2389.700000000001
$ws.Range("N126").Value = -11423

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1514.6471
$ws.Range("I136").Value = 1442.9333
$ws.Range("J136").Value = 2052.5
$ws.Range("K136").Value = 4328.7999
$ws.Range("L136").Value = 6157.5
$ws.Range("M136").Value = -1778.7999
$ws.Range("N136").Value = -11257.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 23638764
$ws.Range("I122").Value = 26002340
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 78007020
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -78004570
$ws.Range("N122").Value = -13900

$ws.Range("H136").Value = 423.52084
$ws.Range("I136").Value = 360.93103
$ws.Range("J136").Value = 519.0526
$ws.Range("K136").Value = 1082.79309
$ws.Range("L136").Value = 1557.1578
$ws.Range("M136").Value = 1467.20691
$ws.Range("N136").Value = -6657.1578
